$wb = $excel.ActiveWorkbook

# --- Model Fit sheet: insert a new "RMSE" row into Table6 (above F-statistic) ---
$ws = $wb.Worksheets.Item("Model Fit")

# Shift rows 7:13 down by inserting a new row at row 7
$ws.Rows.Item(7).Insert()

# Populate the new row with RMSE statistic values
$ws.Range("A7").Value2 = "RMSE"
$ws.Range("B7:D7").NumberFormat = "#,##0.000"
$ws.Range("B7").Value2 = 23.917000000000002
$ws.Range("C7").Value2 = 15.439490694179399
$ws.Range("D7").Value2 = 8.9317796639129803

# Resize the table to include the newly inserted row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D14"))

# Selection left on B13 after the edit, sheet tab active
$ws.Activate()
$ws.Range("B13").Select()

# --- Environmental Model sheet: selection moved ---
$wsEnv = $wb.Worksheets.Item("Environmental Model ")
$wsEnv.Range("A4:A21").Select()

# --- Socioeconomic Model sheet: selection moved, no longer the active tab ---
$wsSoc = $wb.Worksheets.Item("Socioeconomic Model")
$wsSoc.Range("A4:A24").Select()

$ws.Activate()
